# contourmap combined feature engineering
# - when scaling the parameters, the stripe pattern is created
# - there are a lot of differences between the contourmaps of the PCA, but some similarities can be observed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the raw input values for the five Area groups (row 2 = group 0, row 3 = group 1) ---
$ws.Range("B2").Value = 51906
$ws.Range("F2").Value = 25846
$ws.Range("J2").Value = 31791
$ws.Range("N2").Value = 20721
$ws.Range("R2").Value = 50827

$ws.Range("B3").Value = 3507
$ws.Range("F3").Value = 21193
$ws.Range("J3").Value = 22664
$ws.Range("N3").Value = 38758
$ws.Range("R3").Value = 9382

# --- Area4 (N/O columns) percentage formulas get swapped by mistake, creating the "stripe pattern" ---
$ws.Range("O2").Formula = '=N3/$N$4*100'
$ws.Range("O3").Formula = '=N2/$N$4*100'
$ws.Range("O4").Formula = '=N4/$N$4*100'

# --- Totals block (rows 7-10): swap which group is "1" vs "0" ---
$ws.Range("A8").Value = 0
$ws.Range("A9").Value = 1

# B10/C10 recomputed directly from the column totals rather than via the shared formula
$ws.Range("B10").Formula = '=B4+F4+J4+N4+R4'
$ws.Range("C10").Formula = '=B10/$B$10*100'

# --- Selection moves to I9 ---
$ws.Range("I9").Select() | Out-Null

$wb.Save() | Out-Null
